$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($ws, $rowNum, $values)
    $arr = New-Object "object[,]" 1,10
    for ($i = 0; $i -lt 10; $i++) { $arr[0,$i] = $values[$i] }
    $range = $ws.Range("A$rowNum" + ":J$rowNum")
    $range.Value = $arr
}

$ws1 = $wb.Worksheets.Item("Train Results")

Set-RowValues $ws1 2 @(0, 40, 4, 0, 28, 4, 4, 20, 2.31, 2.321960210800171)
Set-RowValues $ws1 3 @(4, 0, 8, 24, 4, 16, 44, 0, 2.98, 3.089115858078003)
Set-RowValues $ws1 4 @(4, 16, 0, 20, 4, 12, 44, 0, 3.22, 3.280884027481079)
Set-RowValues $ws1 5 @(0, 24, 4, 0, 8, 4, 36, 24, 2.45, 2.370848417282104)
Set-RowValues $ws1 6 @(4, 0, 8, 4, 24, 12, 44, 4, 3.17, 3.107202053070068)
Set-RowValues $ws1 7 @(0, 16, 8, 4, 16, 20, 36, 0, 3.07, 2.98982048034668)
Set-RowValues $ws1 8 @(4, 20, 4, 4, 16, 0, 52, 0, 3.35, 3.316265344619751)
Set-RowValues $ws1 9 @(4, 0, 12, 4, 4, 16, 52, 8, 2.81, 2.920425653457642)
Set-RowValues $ws1 10 @(4, 12, 8, 0, 8, 8, 56.00000000000001, 4, 2.88, 2.884921073913574)
Set-RowValues $ws1 11 @(4, 4, 4, 8, 28, 16, 32, 4, 3.04, 3.137231826782227)
Set-RowValues $ws1 12 @(4, 0, 8, 4, 24, 12, 44, 4, 3.09, 3.107202053070068)
Set-RowValues $ws1 13 @(4, 8, 0, 12, 20, 8, 48, 0, 3.11, 3.304723739624023)
Set-RowValues $ws1 14 @(0, 20, 8, 4, 0, 20, 44, 4, 2.82, 2.889007091522217)
Set-RowValues $ws1 15 @(0, 20, 8, 4, 0, 20, 44, 4, 2.74, 2.889007091522217)
Set-RowValues $ws1 16 @(4, 12, 0, 0, 16, 8, 52, 8, 2.62, 2.726036071777344)
Set-RowValues $ws1 17 @(4, 12, 8, 4, 16, 8, 48, 0, 3.19, 3.153226852416992)
Set-RowValues $ws1 18 @(0, 20, 4, 0, 4, 4, 48, 20, 2.52, 2.437030553817749)
Set-RowValues $ws1 19 @(4, 8, 4, 4, 16, 12, 48, 4, 2.94, 3.034233808517456)
Set-RowValues $ws1 20 @(4, 12, 0, 12, 20, 8, 44, 0, 3.29, 3.323781251907349)
Set-RowValues $ws1 21 @(4, 12, 0, 12, 20, 8, 44, 0, 3.25, 3.323781251907349)
Set-RowValues $ws1 22 @(0, 12, 8, 4, 16, 20, 36, 4, 3.01, 2.978699445724487)
Set-RowValues $ws1 23 @(4, 12, 4, 4, 20, 16, 32, 8, 2.96, 3.001132011413574)
Set-RowValues $ws1 24 @(4, 8, 0, 0, 4, 8, 52, 24, 2.44, 2.414699077606201)
Set-RowValues $ws1 25 @(4, 0, 4, 0, 20, 8, 52, 12, 2.92, 2.803555965423584)
Set-RowValues $ws1 26 @(4, 8, 12, 4, 4, 24, 39.99999999999999, 4, 2.9, 2.961641550064087)
Set-RowValues $ws1 27 @(4, 0, 12, 16, 4, 12, 52, 0, 2.96, 3.028928279876709)
Set-RowValues $ws1 28 @(0, 16, 8, 4, 16, 20, 36, 0, 3.02, 2.98982048034668)
Set-RowValues $ws1 29 @(4, 12, 8, 0, 8, 8, 56.00000000000001, 4, 2.79, 2.884921073913574)
Set-RowValues $ws1 30 @(4, 0, 4, 4, 16, 12, 56.00000000000001, 4, 2.85, 2.973666191101074)
Set-RowValues $ws1 31 @(0, 4, 8, 4, 16, 20, 44, 4, 2.94, 2.985975980758667)
Set-RowValues $ws1 32 @(0, 12, 4, 0, 4, 16, 39.99999999999999, 24, 2.51, 2.526475191116333)
Set-RowValues $ws1 33 @(0, 0, 8, 4, 16, 12, 52, 8, 2.99, 2.952923774719238)
Set-RowValues $ws1 34 @(0, 16, 0, 4, 20, 20, 28, 12, 3.38, 2.94190502166748)
Set-RowValues $ws1 35 @(4, 16, 8, 0, 12, 0, 48, 12, 2.56, 2.559431314468384)
Set-RowValues $ws1 36 @(0, 0, 8, 4, 16, 12, 52, 8, 2.82, 2.952923774719238)
Set-RowValues $ws1 37 @(4, 12, 4, 4, 12, 16, 36, 12, 2.86, 2.889096975326538)
Set-RowValues $ws1 38 @(0, 24, 8, 0, 12, 16, 32, 8, 2.93, 2.890259027481079)
Set-RowValues $ws1 39 @(4, 32, 8, 4, 4, 20, 24, 4, 2.84, 2.923932790756226)
Set-RowValues $ws1 40 @(0, 24, 8, 0, 12, 16, 32, 8, 2.94, 2.890259027481079)
Set-RowValues $ws1 41 @(4, 40, 0, 0, 12, 4, 36, 4, 3.16, 3.276581287384033)
Set-RowValues $ws1 42 @(0, 16, 0, 4, 20, 20, 28, 12, 2.72, 2.94190502166748)
Set-RowValues $ws1 43 @(4, 20, 4, 4, 16, 0, 52, 0, 3.21, 3.316265344619751)
Set-RowValues $ws1 44 @(4, 12, 0, 0, 16, 8, 52, 8, 2.56, 2.726036071777344)

$ws2 = $wb.Worksheets.Item("Test Results")

Set-RowValues $ws2 2 @(0, 8, 0, 4, 4, 28, 36, 20, 2.56, 2.822545766830444)
Set-RowValues $ws2 3 @(4, 16, 8, 0, 12, 0, 48, 12, 2.67, 2.559431314468384)
Set-RowValues $ws2 4 @(4, 16, 0, 4, 12, 8, 52, 4, 2.96, 3.057409763336182)
Set-RowValues $ws2 5 @(4, 8, 4, 4, 16, 12, 48, 4, 2.92, 3.034233808517456)
Set-RowValues $ws2 6 @(4, 0, 4, 0, 4, 4, 60, 24, 2.45, 2.419121503829956)
Set-RowValues $ws2 7 @(4, 0, 4, 0, 20, 8, 52, 12, 2.64, 2.803555965423584)
Set-RowValues $ws2 8 @(4, 12, 4, 0, 28, 0, 39.99999999999999, 12, 2.94, 2.711537599563599)
Set-RowValues $ws2 9 @(4, 4, 12, 0, 0, 4, 52, 24, 2.22, 2.446314334869385)
Set-RowValues $ws2 10 @(4, 12, 8, 4, 16, 8, 48, 0, 3.15, 3.153226852416992)
